$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per row, reflecting a repull/recalculation of data.
$updates = @{
    2  = 2
    3  = -3
    4  = 7
    5  = -2
    6  = -3
    7  = 1
    8  = 1
    9  = -4
    11 = 1
    12 = -4
    13 = 6
    14 = -3
    15 = -2
    16 = 2
    17 = 2
    18 = -2
    19 = -4
    20 = 5
    21 = -2
    22 = -1
    23 = 2
    24 = 1
    25 = 5
    26 = -2
    27 = 3
    28 = -1
    29 = -1
    32 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
